$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A248").Value = 37.37
$ws.Range("B248").Value = 0.36
$ws.Range("C248").Value = 46.93020009994507
$ws.Range("D248").Value = 43.59375
$ws.Range("A249").Value = 39.78
$ws.Range("B249").Value = 0.16
$ws.Range("C249").Value = 49.49240255355835
$ws.Range("D249").Value = 43.890625
$ws.Range("A250").Value = 39.78
$ws.Range("B250").Value = 0.16
$ws.Range("C250").Value = 21.56492757797241
$ws.Range("D250").Value = 21.578125
$ws.Range("A251").Value = 39.78
$ws.Range("B251").Value = 0.16
$ws.Range("C251").Value = 26.69873666763306
$ws.Range("D251").Value = 25.03125
$ws.Range("A252").Value = 39.44
$ws.Range("B252").Value = 0.16
$ws.Range("C252").Value = 26.0681004524231
$ws.Range("D252").Value = 24.421875
$ws.Range("A253").Value = 29.44
$ws.Range("B253").Value = 0.16
$ws.Range("C253").Value = 35.41032290458679
$ws.Range("D253").Value = 31.734375
$ws.Range("A254").Value = 100
$ws.Range("B254").Value = 0.16
$ws.Range("C254").Value = 38.95937848091125
$ws.Range("D254").Value = 34.53125
$ws.Range("A255").Value = 23.89
$ws.Range("B255").Value = 0.36
$ws.Range("C255").Value = 68.09541749954224
$ws.Range("D255").Value = 61.390625
$ws.Range("A256").Value = 16.61
$ws.Range("B256").Value = 0.36
$ws.Range("C256").Value = 43.96134662628174
$ws.Range("D256").Value = 40.328125
$ws.Range("A257").Value = 19.12
$ws.Range("B257").Value = 0.16
$ws.Range("C257").Value = 34.96498870849609
$ws.Range("D257").Value = 33.15625
$ws.Range("A258").Value = 75
$ws.Range("B258").Value = 0.0004
$ws.Range("C258").Value = 1.809213399887085
$ws.Range("D258").Value = 1.5
$ws.Range("A259").Value = 75
$ws.Range("B259").Value = 0.0004
$ws.Range("C259").Value = 1.765857219696045
$ws.Range("D259").Value = 1.375
$ws.Range("A260").Value = 39.44
$ws.Range("B260").Value = 0.16
$ws.Range("C260").Value = 41.00721597671509
$ws.Range("D260").Value = 39.75
$ws.Range("A261").Value = 37.67
$ws.Range("B261").Value = 0.36
$ws.Range("C261").Value = 58.08780217170715
$ws.Range("D261").Value = 54.03125
$ws.Range("A262").Value = 37.67
$ws.Range("B262").Value = 0.36
$ws.Range("C262").Value = 69.58404088020325
$ws.Range("D262").Value = 67.546875
$ws.Range("A263").Value = 39.44
$ws.Range("B263").Value = 0.16
$ws.Range("C263").Value = 47.34826254844666
$ws.Range("D263").Value = 45.734375
$ws.Range("A264").Value = 39.44
$ws.Range("B264").Value = 0.16
$ws.Range("C264").Value = 49.34514284133911
$ws.Range("D264").Value = 47.875
$ws.Range("A265").Value = 39.44
$ws.Range("B265").Value = 0.16
$ws.Range("C265").Value = 44.41006135940552
$ws.Range("D265").Value = 42.953125
$ws.Range("A266").Value = 37
$ws.Range("B266").Value = 0.01
$ws.Range("C266").Value = 3.879632472991943
$ws.Range("D266").Value = 2.890625
$ws.Range("A267").Value = 75
$ws.Range("B267").Value = 0.0004
$ws.Range("C267").Value = 2.40252685546875
$ws.Range("D267").Value = 2.09375
$ws.Range("A268").Value = 50
$ws.Range("B268").Value = 0.0016
$ws.Range("C268").Value = 2.566594839096069
$ws.Range("D268").Value = 2.09375
$ws.Range("A269").Value = 39.44
$ws.Range("B269").Value = 0.16
$ws.Range("C269").Value = 36.68801498413086
$ws.Range("D269").Value = 35.109375
$ws.Range("A270").Value = 37.67
$ws.Range("B270").Value = 0.36
$ws.Range("C270").Value = 52.79992604255676
$ws.Range("D270").Value = 51.234375
$ws.Range("A271").Value = 75
$ws.Range("B271").Value = 0.0004
$ws.Range("C271").Value = 1.868182420730591
$ws.Range("D271").Value = 1.390625
